# Update Efnb2-Ephb4 LR-pair sheet with recomputed TPM-based values.
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.76217133333333
$ws.Range("H2").Value = 137.286514
$ws.Range("I2").Value = 0.6763939203605134
$ws.Range("J2").Value = 0.6763939203605135
$ws.Range("M2").Value = 34.82273866666667
$ws.Range("N2").Value = 104.468216
$ws.Range("O2").Value = 0.7026863693863559
$ws.Range("P2").Value = 0.702686369386356
$ws.Range("Q2").Value = 1593.564133159892
$ws.Range("R2").Value = 14342.07719843902
$ws.Range("S2").Value = 0.4752927881731331
$ws.Range("T2").Value = 0.4752927881731333

# Row 3
$ws.Range("G3").Value = 45.76217133333333
$ws.Range("H3").Value = 137.286514
$ws.Range("I3").Value = 0.6763939203605134
$ws.Range("J3").Value = 0.6763939203605135
$ws.Range("O3").Value = 0.1722322446965897
$ws.Range("P3").Value = 0.1722322446965897
$ws.Range("Q3").Value = 390.5912220295198
$ws.Range("R3").Value = 3515.320998265679
$ws.Range("S3").Value = 0.1164968432028176
$ws.Range("T3").Value = 0.1164968432028176

# Row 4
$ws.Range("G4").Value = 45.76217133333333
$ws.Range("H4").Value = 137.286514
$ws.Range("I4").Value = 0.6763939203605134
$ws.Range("J4").Value = 0.6763939203605135
$ws.Range("M4").Value = 6.169951999999999
$ws.Range("N4").Value = 18.509856
$ws.Range("O4").Value = 0.1245031647760143
$ws.Range("P4").Value = 0.1245031647760143
$ws.Range("Q4").Value = 282.3504005424426
$ws.Range("R4").Value = 2541.153604881984
$ws.Range("S4").Value = 0.08421318372013929
$ws.Range("T4").Value = 0.08421318372013933

# Row 5
$ws.Range("G5").Value = 45.76217133333333
$ws.Range("H5").Value = 137.286514
$ws.Range("I5").Value = 0.6763939203605134
$ws.Range("J5").Value = 0.6763939203605135
$ws.Range("M5").Value = 0.02865466666666667
$ws.Range("N5").Value = 0.085964
$ws.Range("O5").Value = 0.000578221141039957
$ws.Range("P5").Value = 0.0005782211410399571
$ws.Range("Q5").Value = 1.311299765499556
$ws.Range("R5").Value = 11.801697889496
$ws.Range("S5").Value = 0.0003911052644233458
$ws.Range("T5").Value = 0.000391105264423346

# Row 6
$ws.Range("I6").Value = 0.1388778842960613
$ws.Range("J6").Value = 0.1388778842960613
$ws.Range("M6").Value = 34.82273866666667
$ws.Range("N6").Value = 104.468216
$ws.Range("O6").Value = 0.7026863693863559
$ws.Range("P6").Value = 0.702686369386356
$ws.Range("Q6").Value = 327.1922006415662
$ws.Range("R6").Value = 2944.729805774096
$ws.Range("S6").Value = 0.09758759630405772
$ws.Range("T6").Value = 0.09758759630405775

# Row 7
$ws.Range("I7").Value = 0.1388778842960613
$ws.Range("J7").Value = 0.1388778842960613
$ws.Range("O7").Value = 0.1722322446965897
$ws.Range("P7").Value = 0.1722322446965897
$ws.Range("S7").Value = 0.0239192497510239
$ws.Range("T7").Value = 0.02391924975102391

# Row 8
$ws.Range("I8").Value = 0.1388778842960613
$ws.Range("J8").Value = 0.1388778842960613
$ws.Range("M8").Value = 6.169951999999999
$ws.Range("N8").Value = 18.509856
$ws.Range("O8").Value = 0.1245031647760143
$ws.Range("P8").Value = 0.1245031647760143
$ws.Range("Q8").Value = 57.97247000177065
$ws.Range("R8").Value = 521.752230015936
$ws.Range("S8").Value = 0.01729073611225677
$ws.Range("T8").Value = 0.01729073611225678

# Row 9
$ws.Range("I9").Value = 0.1388778842960613
$ws.Range("J9").Value = 0.1388778842960613
$ws.Range("M9").Value = 0.02865466666666667
$ws.Range("N9").Value = 0.085964
$ws.Range("O9").Value = 0.000578221141039957
$ws.Range("P9").Value = 0.0005782211410399571
$ws.Range("Q9").Value = 0.2692373949982222
$ws.Range("R9").Value = 2.423136554984
$ws.Range("S9").Value = 0.00008030212872288369
$ws.Range("T9").Value = 0.00008030212872288372

# Row 10
$ws.Range("G10").Value = 12.29750866666667
$ws.Range("H10").Value = 36.892526
$ws.Range("I10").Value = 0.1817649787009828
$ws.Range("J10").Value = 0.1817649787009828
$ws.Range("M10").Value = 34.82273866666667
$ws.Range("N10").Value = 104.468216
$ws.Range("O10").Value = 0.7026863693863559
$ws.Range("P10").Value = 0.702686369386356
$ws.Range("Q10").Value = 428.2329305504018
$ws.Range("R10").Value = 3854.096374953616
$ws.Range("S10").Value = 0.1277237729649819
$ws.Range("T10").Value = 0.1277237729649819

# Row 11
$ws.Range("G11").Value = 12.29750866666667
$ws.Range("H11").Value = 36.892526
$ws.Range("I11").Value = 0.1817649787009828
$ws.Range("J11").Value = 0.1817649787009828
$ws.Range("O11").Value = 0.1722322446965897
$ws.Range("P11").Value = 0.1722322446965897
$ws.Range("Q11").Value = 104.9622165662669
$ws.Range("R11").Value = 944.6599490964021
$ws.Range("S11").Value = 0.03130579028889808
$ws.Range("T11").Value = 0.03130579028889809

# Row 12
$ws.Range("G12").Value = 12.29750866666667
$ws.Range("H12").Value = 36.892526
$ws.Range("I12").Value = 0.1817649787009828
$ws.Range("J12").Value = 0.1817649787009828
$ws.Range("M12").Value = 6.169951999999999
$ws.Range("N12").Value = 18.509856
$ws.Range("O12").Value = 0.1245031647760143
$ws.Range("P12").Value = 0.1245031647760143
$ws.Range("Q12").Value = 75.87503819291733
$ws.Range("R12").Value = 682.875343736256
$ws.Range("S12").Value = 0.02263031509371719
$ws.Range("T12").Value = 0.0226303150937172

# Row 13
$ws.Range("G13").Value = 12.29750866666667
$ws.Range("H13").Value = 36.892526
$ws.Range("I13").Value = 0.1817649787009828
$ws.Range("J13").Value = 0.1817649787009828
$ws.Range("M13").Value = 0.02865466666666667
$ws.Range("N13").Value = 0.085964
$ws.Range("O13").Value = 0.000578221141039957
$ws.Range("P13").Value = 0.0005782211410399571
$ws.Range("Q13").Value = 0.3523810116737778
$ws.Range("R13").Value = 3.171429105064
$ws.Range("S13").Value = 0.0001051003533855857
$ws.Range("T13").Value = 0.0001051003533855858

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2004796666666666
$ws.Range("H14").Value = 0.6014389999999999
$ws.Range("I14").Value = 0.002963216642442438
$ws.Range("J14").Value = 0.002963216642442439
$ws.Range("M14").Value = 34.82273866666667
$ws.Range("N14").Value = 104.468216
$ws.Range("O14").Value = 0.7026863693863559
$ws.Range("P14").Value = 0.702686369386356
$ws.Range("Q14").Value = 6.981251040313777
$ws.Range("R14").Value = 62.83125936282399
$ws.Range("S14").Value = 0.002082211944183104
$ws.Range("T14").Value = 0.002082211944183105

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2004796666666666
$ws.Range("H15").Value = 0.6014389999999999
$ws.Range("I15").Value = 0.002963216642442438
$ws.Range("J15").Value = 0.002963216642442439
$ws.Range("O15").Value = 0.1722322446965897
$ws.Range("P15").Value = 0.1722322446965897
$ws.Range("Q15").Value = 1.711142537905889
$ws.Range("R15").Value = 15.400282841153
$ws.Range("S15").Value = 0.000510361453850153
$ws.Range("T15").Value = 0.0005103614538501532

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2004796666666666
$ws.Range("H16").Value = 0.6014389999999999
$ws.Range("I16").Value = 0.002963216642442438
$ws.Range("J16").Value = 0.002963216642442439
$ws.Range("M16").Value = 6.169951999999999
$ws.Range("N16").Value = 18.509856
$ws.Range("O16").Value = 0.1245031647760143
$ws.Range("P16").Value = 0.1245031647760143
$ws.Range("Q16").Value = 1.236949920309333
$ws.Range("R16").Value = 11.132549282784
$ws.Range("S16").Value = 0.0003689298499010388
$ws.Range("T16").Value = 0.0003689298499010389

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2004796666666666
$ws.Range("H17").Value = 0.6014389999999999
$ws.Range("I17").Value = 0.002963216642442438
$ws.Range("J17").Value = 0.002963216642442439
$ws.Range("M17").Value = 0.02865466666666667
$ws.Range("N17").Value = 0.085964
$ws.Range("O17").Value = 0.000578221141039957
$ws.Range("P17").Value = 0.0005782211410399571
$ws.Range("Q17").Value = 0.005744678021777777
$ws.Range("R17").Value = 0.05170210219599999
$ws.Range("S17").Value = 0.000001713394508141657
$ws.Range("T17").Value = 0.000001713394508141658
